$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply the date number format to the new column's cells first so the
# format (and derived styles) get created before any values are written.
$ws.Range("C1:C4").NumberFormat = "dd.mm.yyyy"

# Header for the new "date_time" column.
$ws.Range("C1").Value = "date_time"

# Data rows: 2020-01-01 for every row (matches serial date 43831).
$d = Get-Date -Year 2020 -Month 1 -Day 1 -Hour 0 -Minute 0 -Second 0
$ws.Range("C2").Value = $d
$ws.Range("C3").Value = $d
$ws.Range("C4").Value = $d

# Header cell is right-aligned (matches the added cellXfs entry).
$ws.Range("C1").HorizontalAlignment = -4152

# Match the workbook's final selection/active cell.
$ws.Range("C4").Select()
